$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update / add rows 2-10 with the Alcam-Chl1 ECs/FAPs/sCs ligand-receptor matrix
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Alcam"
$ws.Range("C2").Value = "Chl1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 62.12558000000001
$ws.Range("H2").Value = 186.37674
$ws.Range("I2").Value = 0.9736910227596813
$ws.Range("J2").Value = 0.9736910227596813
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.041336
$ws.Range("N2").Value = 0.124008
$ws.Range("O2").Value = 0.01439314688988224
$ws.Range("P2").Value = 0.01439314688988224
$ws.Range("Q2").Value = 2.568022974880001
$ws.Range("R2").Value = 23.11220677392
$ws.Range("S2").Value = 0.01401447791593976
$ws.Range("T2").Value = 0.01401447791593976

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Alcam"
$ws.Range("C3").Value = "Chl1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 62.12558000000001
$ws.Range("H3").Value = 186.37674
$ws.Range("I3").Value = 0.9736910227596813
$ws.Range("J3").Value = 0.9736910227596813
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.466778
$ws.Range("N3").Value = 1.400334
$ws.Range("O3").Value = 0.1625315540682565
$ws.Range("P3").Value = 0.1625315540682565
$ws.Range("Q3").Value = 28.99885398124
$ws.Range("R3").Value = 260.98968583116
$ws.Range("S3").Value = 0.1582555151114411
$ws.Range("T3").Value = 0.1582555151114411

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Alcam"
$ws.Range("C4").Value = "Chl1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 62.12558000000001
$ws.Range("H4").Value = 186.37674
$ws.Range("I4").Value = 0.9736910227596813
$ws.Range("J4").Value = 0.9736910227596813
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.363808333333333
$ws.Range("N4").Value = 7.091424999999999
$ws.Range("O4").Value = 0.8230752990418613
$ws.Range("P4").Value = 0.8230752990418614
$ws.Range("Q4").Value = 146.8529637171667
$ws.Range("R4").Value = 1321.6766734545
$ws.Range("S4").Value = 0.8014210297323004
$ws.Range("T4").Value = 0.8014210297323006

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Alcam"
$ws.Range("C5").Value = "Chl1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5683613333333334
$ws.Range("H5").Value = 1.705084
$ws.Range("I5").Value = 0.008907897969731461
$ws.Range("J5").Value = 0.008907897969731461
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.041336
$ws.Range("N5").Value = 0.124008
$ws.Range("O5").Value = 0.01439314688988224
$ws.Range("P5").Value = 0.01439314688988224
$ws.Range("Q5").Value = 0.02349378407466667
$ws.Range("R5").Value = 0.211444056672
$ws.Range("S5").Value = 0.0001282126839584287
$ws.Range("T5").Value = 0.0001282126839584287

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Alcam"
$ws.Range("C6").Value = "Chl1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5683613333333334
$ws.Range("H6").Value = 1.705084
$ws.Range("I6").Value = 0.008907897969731461
$ws.Range("J6").Value = 0.008907897969731461
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.466778
$ws.Range("N6").Value = 1.400334
$ws.Range("O6").Value = 0.1625315540682565
$ws.Range("P6").Value = 0.1625315540682565
$ws.Range("Q6").Value = 0.2652985664506667
$ws.Range("R6").Value = 2.387687098056
$ws.Range("S6").Value = 0.001447814500501921
$ws.Range("T6").Value = 0.001447814500501921

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Alcam"
$ws.Range("C7").Value = "Chl1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5683613333333334
$ws.Range("H7").Value = 1.705084
$ws.Range("I7").Value = 0.008907897969731461
$ws.Range("J7").Value = 0.008907897969731461
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.363808333333333
$ws.Range("N7").Value = 7.091424999999999
$ws.Range("O7").Value = 0.8230752990418613
$ws.Range("P7").Value = 0.8230752990418614
$ws.Range("Q7").Value = 1.343497256077778
$ws.Range("R7").Value = 12.0914753047
$ws.Range("S7").Value = 0.007331870785271112
$ws.Range("T7").Value = 0.007331870785271113

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Alcam"
$ws.Range("C8").Value = "Chl1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.110262
$ws.Range("H8").Value = 3.330786
$ws.Range("I8").Value = 0.01740107927058724
$ws.Range("J8").Value = 0.01740107927058724
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.041336
$ws.Range("N8").Value = 0.124008
$ws.Range("O8").Value = 0.01439314688988224
$ws.Range("P8").Value = 0.01439314688988224
$ws.Range("Q8").Value = 0.04589379003200001
$ws.Range("R8").Value = 0.413044110288
$ws.Range("S8").Value = 0.000250456289984047
$ws.Range("T8").Value = 0.0002504562899840471

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Alcam"
$ws.Range("C9").Value = "Chl1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.110262
$ws.Range("H9").Value = 3.330786
$ws.Range("I9").Value = 0.01740107927058724
$ws.Range("J9").Value = 0.01740107927058724
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.466778
$ws.Range("N9").Value = 1.400334
$ws.Range("O9").Value = 0.1625315540682565
$ws.Range("P9").Value = 0.1625315540682565
$ws.Range("Q9").Value = 0.518245875836
$ws.Range("R9").Value = 4.664212882524001
$ws.Range("S9").Value = 0.002828224456313467
$ws.Range("T9").Value = 0.002828224456313468

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Alcam"
$ws.Range("C10").Value = "Chl1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.110262
$ws.Range("H10").Value = 3.330786
$ws.Range("I10").Value = 0.01740107927058724
$ws.Range("J10").Value = 0.01740107927058724
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.363808333333333
$ws.Range("N10").Value = 7.091424999999999
$ws.Range("O10").Value = 0.8230752990418613
$ws.Range("P10").Value = 0.8230752990418614
$ws.Range("Q10").Value = 2.624446567783333
$ws.Range("R10").Value = 23.62001911005
$ws.Range("S10").Value = 0.01432239852428973
$ws.Range("T10").Value = 0.01432239852428973

